$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row holding only the professor name (B13/C13, no label in A13) is
# removed entirely; everything below shifts up one row (row heights follow).
$ws.Rows("13").Delete()

# Row 10 (Objetivos:) now shows the professor name instead of the long
# Portuguese objectives paragraph.
$ws.Range("B10").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C10").Value = "8855158 - Morun Bernardino Neto"

# Row 13 (Programa resumido:) now shows "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date "01/01/2022" (same text
# already used by A8/B8). Use copy/paste-values instead of a plain string
# assignment so the date-like text stays plain text instead of being
# reinterpreted as a date serial number.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 18 (Método:) now shows the professor name.
$ws.Range("B18").Value = "8855158 - Morun Bernardino Neto"
$ws.Range("C18").Value = "8855158 - Morun Bernardino Neto"

# Row 19 (Critério:) now shows the evaluation method text.
$ws.Range("B19").Value = "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais."
$ws.Range("C19").Value = "O método de avaliação será composto por avaliações teóricas, práticas e relatórios de atividades de práticas laboratoriais."

# Row 20 (Norma de recuperação:) now shows the final grade criterion text.
$ws.Range("B20").Value = "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."
$ws.Range("C20").Value = "Para o cálculo da nota final (NF) será feita a média aritmética das avaliações aplicadas. Estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."

# Row 21 (Bibliografia:) now shows the recovery evaluation norm text.
$ws.Range("B21").Value = "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."
$ws.Range("C21").Value = "Avaliação de recuperação (R) envolvendo todo o conteúdo da disciplina. Média Final = (NF+R) / 2 e estará aprovado por notas o aluno que obtiver nota final igual ou superior a 5,0 pontos."
